$wb = $excel.ActiveWorkbook

# Sheet "OFF" - Week 17 row (row 2) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 278
$wsOff.Range("C2").Value = 172
$wsOff.Range("D2").Value = 55
$wsOff.Range("E2").Value = 29

# Sheet "DEF" - Week 17 row (row 2) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 264
$wsDef.Range("C2").Value = 198
$wsDef.Range("D2").Value = 74
$wsDef.Range("E2").Value = 45
$wsDef.Range("G2").Value = 3
